# Kenntnisse.xlsx update:
#  - add two new IT skills ("Power Auomate", "Power Apps") right before "SAC / SAP"
#  - drop the "Französisch / B1" language row
#  - widen column A and move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows just above the "SAC / SAP" row (row 8), pushing every
# row below it (SAC/SAP, Deutsch, Englisch, Französisch) down by two.
$ws.Rows("8:9").Insert()

# Row 8: Power Auomate
$ws.Range("A8").Value = "Power Auomate"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = "IT"

# Row 9: Power Apps
$ws.Range("A9").Value = "Power Apps"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "IT"

# The "Französisch / B1" language entry (now shifted down to row 13) is no
# longer wanted - remove the whole row.
$ws.Rows(13).Delete()

# Column A now holds longer labels ("Power Auomate" / "Power Apps") - widen it.
$ws.Columns(1).ColumnWidth = 15.5

# Leave the selection where the author left it after editing the table.
$ws.Range("A14").Select() | Out-Null
